# Insert a new weekly price-report row for "Feria Lagunitas de Puerto Montt -
# Ciboulette" above the current row 100, pushing the existing rows 100-118
# down to 101-119 (dimension grows from A1:R118 to A1:R119).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row at row 100 (shifts rows 100..118 down to 101..119).
$ws.Rows.Item(100).EntireRow.Insert()

# Populate the newly inserted row 100 with the new data record.
$ws.Cells.Item(100, 1).Value = 4
$ws.Cells.Item(100, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(100, 3).Value = "Los Lagos"
$ws.Cells.Item(100, 4).Value = 44476
$ws.Cells.Item(100, 5).Value = 10
$ws.Cells.Item(100, 6).Value = 100112039
$ws.Cells.Item(100, 7).Value = "Ciboulette"
$ws.Cells.Item(100, 8).Value = "Sin especificar"
$ws.Cells.Item(100, 9).Value = "Primera"
$ws.Cells.Item(100, 10).Value = 120
$ws.Cells.Item(100, 11).Value = 3000
$ws.Cells.Item(100, 12).Value = 3000
$ws.Cells.Item(100, 13).Value = 3000
$ws.Cells.Item(100, 14).Value = "$/docena de atados"
$ws.Cells.Item(100, 15).Value = "Región Metropolitana"
$ws.Cells.Item(100, 16).Value = 1000
$ws.Cells.Item(100, 17).Value = 3
$ws.Cells.Item(100, 18).Value = "Hortaliza"
